$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "protocol"
$ws.Range("L2").Value = "GO-SHIP"
$ws.Range("M1").Value = "responsible_group"
$ws.Range("M2").Value = "Davis Strait Observing System"
$ws.Range("M3").Value = "RAPID"
$ws.Range("L4").Value = "GO-SHIP"
$ws.Range("M4").Value = "Davis Strait Observing System"

$ws.Range("M4").Select()
